$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update threshold values (rows 2-4 keep their parameter name, only
#     Min/Max numbers change) ---
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 12

$ws.Range("B3").Value = 4.0999999999999996
$ws.Range("C3").Value = 10

$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.4

# --- Remove the "theta_threshold_range" row entirely (old row 5). This
#     shifts "pie_threshold_range" (old row 6) up into row 5, keeping its
#     own Min value (0, unchanged) and Max value needs updating to 20.
#     Deleting the row (rather than overwriting text) also drops the now-
#     unused "theta_threshold_range" shared string, matching the trimmed
#     sharedStrings table in the target file. ---
$ws.Rows("5").Delete()

$ws.Range("C5").Value = 20

# --- Page setup: portrait orientation, A4 paper ---
$ps = $ws.PageSetup
$ps.Orientation = 1
$ps.PaperSize = 9

# --- Selection moves to C3 ---
$ws.Range("C3").Select() | Out-Null
